$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 additions
$ws.Range("X7").Value = 0.47000100000001055
$ws.Range("Y7").Value = "Up"

# Row 8 - new row
$ws.Range("A8").Value = 42649.886817129627
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = 21
$ws.Range("C8").Value = "Strong Buy"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = "Random"
$ws.Range("Q8").Value = 35.958706302092025
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = -0.0019
$ws.Range("T8").Value = -0.023
$ws.Range("S7:T7").Copy()
$ws.Range("S8:T8").PasteSpecial(-4122)
$ws.Range("U8").Value = 14.71
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0

# Column C width change (closest achievable width to the target 9.125 via the
# COM ColumnWidth setter, which snaps to whole-pixel granularity)
$ws.Columns("C").ColumnWidth = 8.29
